$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure target cells keep a plain "Text" number format so values
# like "44.015.28" or "0.0760" are stored verbatim as strings,
# matching the inlineStr cells produced by the original export.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '44.015.28'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +0.27%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.357.02'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +0.14%  '

$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.25%  '

$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +0.97%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '239.45'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +0.87%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '74.25'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +2.23%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.588'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +8.64%  '

$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +1.07%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '57.19'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +0.01%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '31.98'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +13.56%  '

$ws.Range('B13').NumberFormat = '@'
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').NumberFormat = '@'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '7.25'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +8.93%  '

$ws.Range('B14').NumberFormat = '@'
$ws.Range('B14').Value = 'TRON'
$ws.Range('C14').NumberFormat = '@'
$ws.Range('C14').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.107'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +0.39%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.706.97'

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '16.62'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -0.30%  '

$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +1.01%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.365.73'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +0.27%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '43.891.31'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +0.02%  '

$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +0.46%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.78'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +5.17%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '76.89'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -1.12%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '256.06'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +0.67%  '

$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +23.35%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.999'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -0.05%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.69'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -1.73%  '

$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -0.57%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.70'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +1.46%  '

$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -2.32%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '22.73'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +1.47%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '175.52'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +1.86%  '

$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -2.30%  '

$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +3.81%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0760'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +6.76%  '

$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +1.56%  '

$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +3.94%  '

$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -6.56%  '

$ws.Range('B38').NumberFormat = '@'
$ws.Range('B38').Value = 'LidoDAOToken'
$ws.Range('C38').NumberFormat = '@'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.36'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -2.84%  '

$ws.Range('B39').NumberFormat = '@'
$ws.Range('B39').Value = 'THORChain'
$ws.Range('C39').NumberFormat = '@'
$ws.Range('C39').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '6.32'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -1.57%  '

$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +4.09%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.114'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +16.74%  '

$ws.Range('B42').NumberFormat = '@'
$ws.Range('B42').Value = 'InjectiveProtocol'
$ws.Range('C42').NumberFormat = '@'
$ws.Range('C42').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '19.25'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -2.29%  '

$ws.Range('B43').NumberFormat = '@'
$ws.Range('B43').Value = 'FraxShare'
$ws.Range('C43').NumberFormat = '@'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '9.16'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +3.74%  '

$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +11.32%  '

$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -0.01%  '

$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +5.82%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '57.89'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +10.30%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.48'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +7.61%  '

$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +0.24%  '

$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +0.93%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '100.23'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +2.75%  '
